$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.355.25'
$ws.Range("E2").Value = '  -2.06%  '

$ws.Range("D3").Value = '3.494.11'
$ws.Range("E3").Value = '  -2.05%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '609.91'
$ws.Range("E5").Value = '  +4.63%  '

$ws.Range("D6").Value = '186.47'
$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("D7").Value = '0.627'
$ws.Range("E7").Value = '  +0.98%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").Value = '0.212'
$ws.Range("E9").Value = '  -3.58%  '

$ws.Range("D10").Value = '0.653'
$ws.Range("E10").Value = '  +0.03%  '

$ws.Range("D11").Value = '53.23'
$ws.Range("E11").Value = '  -2.69%  '

$ws.Range("E12").Value = '  -3.95%  '

$ws.Range("D13").Value = '9.62'
$ws.Range("E13").Value = '  +0.91%  '

$ws.Range("D14").Value = '4.051.27'
$ws.Range("E14").Value = '  +0.36%  '

$ws.Range("D15").Value = '605.41'
$ws.Range("E15").Value = '  +7.14%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '69.462.01'
$ws.Range("E16").Value = '  -1.89%  '

$ws.Range("D17").Value = '12.68'
$ws.Range("E17").Value = '  +1.75%  '

$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = '18.92'
$ws.Range("E18").Value = '  -2.21%  '

$ws.Range("D19").Value = '3.484.12'
$ws.Range("E19").Value = '  -2.47%  '

$ws.Range("E20").Value = '  -0.37%  '

$ws.Range("E21").Value = '  -1.98%  '

$ws.Range("D22").Value = '17.35'
$ws.Range("E22").Value = '  -2.54%  '

$ws.Range("D23").Value = '104.77'
$ws.Range("E23").Value = '  +11.20%  '

$ws.Range("D24").Value = '4.65'
$ws.Range("E24").Value = '  +1.47%  '

$ws.Range("D25").Value = '5.02'
$ws.Range("E25").Value = '  -0.44%  '

$ws.Range("D26").Value = '3.05'
$ws.Range("E26").Value = '  +3.36%  '

$ws.Range("D27").Value = '10.94'
$ws.Range("E27").Value = '  -3.37%  '

$ws.Range("D28").Value = '9.91'
$ws.Range("E28").Value = '  +7.90%  '

$ws.Range("D29").Value = '33.73'
$ws.Range("E29").Value = '  +3.73%  '

$ws.Range("D30").Value = '7.01'
$ws.Range("E30").Value = '  -3.88%  '

$ws.Range("E31").Value = '  +1.27%  '

$ws.Range("E32").Value = '  -0.31%  '

$ws.Range("B33").Value = 'dogwifhat'
$ws.Range("C33").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D33").Value = '3.86'
$ws.Range("E33").Value = '  +14.47%  '

$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").Value = '63.33'
$ws.Range("E34").Value = '  -0.50%  '

$ws.Range("E35").Value = '  -6.77%  '

$ws.Range("E36").Value = '  -0.11%  '

$ws.Range("D37").Value = '524.01'
$ws.Range("E37").Value = '  -5.38%  '

$ws.Range("E38").Value = '  -5.57%  '

$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '3.568.28'
$ws.Range("E39").Value = '  +0.36%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '3.58'
$ws.Range("E40").Value = '  +3.62%  '

$ws.Range("D41").Value = '36.68'
$ws.Range("E41").Value = '  -2.88%  '

$ws.Range("D42").Value = '0.0₃0770'
$ws.Range("E42").Value = '  -4.19%  '

$ws.Range("D43").Value = '0.138'
$ws.Range("E43").Value = '  +1.16%  '

$ws.Range("D44").Value = '0.0462'
$ws.Range("E44").Value = '  +2.68%  '

$ws.Range("E45").Value = '  +1.27%  '

$ws.Range("D46").Value = '0.145'
$ws.Range("E46").Value = '  +5.77%  '

$ws.Range("D47").Value = '3.34'
$ws.Range("E47").Value = '  -4.67%  '

$ws.Range("E48").Value = '  -5.30%  '

$ws.Range("E49").Value = '  +0.42%  '

$ws.Range("D50").Value = '131.14'
$ws.Range("E50").Value = '  -3.41%  '

$ws.Range("E51").Value = '  -8.73%  '
